$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix capitalization of "luke" -> "Luke" (Person column) ---
$ws.Range("F22").Value = "Luke"

# --- Row 26: rework the "mouse" task row ---
# Shorten text and mark done, add progress note, fix person capitalization
$ws.Range("A26").Value = "get mouse to work"
$ws.Range("B26").Value = "done"
$ws.Range("D26").Value = "could be more efficient"
$ws.Range("F26").Value = "Luke"

# --- Fix spelling: "Asthetic:" -> "Aesthetic:" ---
$ws.Range("A27").Value = "Aesthetic:"

# --- Row 32: "Obscure vision" gamefeel task gains status/notes/person ---
$ws.Range("B32").Value = "playable"
$ws.Range("D32").Value = "worked too well - finetune"
$ws.Range("F32").Value = "luke"
$ws.Rows(32).RowHeight = 30

# --- Insert two new rows before the old last "Artist Polish" row, for AI section ---
$ws.Rows(34).Insert()
$ws.Rows(34).Insert()

# Copy formatting (border style) from a neighboring plain row onto the new rows
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A34:A35").PasteSpecial(-4122) | Out-Null

$ws.Range("A34").Value = "AI Behavior"
$ws.Range("A35").Value = "Gameplay objects"

# --- Update selection / view to reflect the edited area ---
$ws.Range("D34").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
